$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

# Date header
Replace-Text "2024-01-19 Friday" "2024-01-20 Saturday"

# Table cell values.
# NOTE: "53÷3=" is both the target of one replacement (66÷7= -> 53÷3=)
# and the source of another (53÷3= -> 80÷7=). Do the 53÷3=->80÷7=
# replacement first so that we don't re-match the freshly-written text.
Replace-Text "53÷3=" "80÷7="

Replace-Text "66÷7=" "53÷3="
Replace-Text "64÷6=" "95÷8="
Replace-Text "99÷8=" "27÷2="
Replace-Text "19÷4=" "10÷3="
Replace-Text "58÷3=" "93÷8="
Replace-Text "64÷3=" "92÷8="
Replace-Text "29÷7=" "35÷5="
Replace-Text "45÷8=" "34÷2="
Replace-Text "40÷3=" "93÷9="
Replace-Text "17÷4=" "37÷5="
Replace-Text "96÷9=" "39÷6="
Replace-Text "54÷2=" "56÷9="
Replace-Text "15÷5=" "55÷6="
Replace-Text "54÷9=" "91÷7="
Replace-Text "82÷3=" "42÷5="
Replace-Text "17÷6=" "20÷4="
Replace-Text "54÷3=" "87÷4="
Replace-Text "66÷2=" "65÷6="
Replace-Text "70÷2=" "98÷3="
Replace-Text "25÷9=" "83÷9="
Replace-Text "68÷4=" "20÷3="
Replace-Text "18÷6=" "94÷6="
Replace-Text "94÷3=" "87÷4="
Replace-Text "17÷5=" "42÷5="
